$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210238695144653
$ws.Range("B1").Value = 1.435527801513672
$ws.Range("C1").Value = 6.898724555969238
$ws.Range("D1").Value = 2.18938422203064
$ws.Range("E1").Value = 1.17006778717041
